# "Changes of 22nd June 2022"
# The route work's ready/scheduled dates move from 2022-06-17 to 2022-06-22
# (serial dates 44729 -> 44734) and the associated times move from
# 07:00/06:45 to 22:00/21:45 (i.e. +15 hours), mirrored in both the
# "primary" columns (P/Q/R/S/T) and the duplicate tracking columns
# (CA/CB/CC).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Ready date/time
$ws.Range("P2").Value = 44734
$ws.Range("S2").Value = 0.91666666666666663
$ws.Range("T2").Value = 0.90625

# Duplicate tracking columns (same date/time as P2/S2)
$ws.Range("CA2").Value = 44734
$ws.Range("CB2").Value = 0.91666666666666663

# Reflect the scroll/selection state the author had when saving: the
# workbook view was scrolled back to column M and P2 (the date that was
# just edited) was the active cell.
$ws.Activate()
$ws.Range("P2").Select()
